# Commit Upto Sales Person
#
# 1) Re-work a few rows of the "SalesTypes" sheet (new sales-type values,
#    a hyperlink on the type code, and a couple of value tweaks).
# 2) Add a brand-new "SalesPersonData" sheet listing sales people and
#    their provision percentage.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "SalesTypes" sheet edits
# ---------------------------------------------------------------------
$salesTypes = $wb.Worksheets.Item("SalesTypes")

$salesTypes.Range("A2").Value = "Transport"
$salesTypes.Range("B2").Value = 1.33

$salesTypes.Range("A3").Value = "C@B"
$salesTypes.Range("B3").Value = 6.5
$salesTypes.Hyperlinks.Add($salesTypes.Range("A3"), "http://www.bgdad.com")

$salesTypes.Range("A4").Value = 56565
$salesTypes.Range("B4").Value = "B6525D"

[void]$salesTypes.Range("B10:B11").Select()
$null = $salesTypes.Range("B11").Activate()

# ---------------------------------------------------------------------
# 2) New "SalesPersonData" sheet, appended after the last tab
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$salesPerson = $wb.Worksheets.Add($null, $lastSheet)
$salesPerson.Name = "SalesPersonData"

$salesPerson.Columns.Item(1).ColumnWidth = 17.85546875

$salesPerson.Range("A1").Value = "Sales person name"
$salesPerson.Range("B1").Value = "Provision"
$salesPerson.Range("B1").Font.Name = "Verdana"
$salesPerson.Range("B1").Font.Size = 8
$salesPerson.Range("B1").Font.Color = 0
$salesPerson.Range("B1").HorizontalAlignment = -4108

$salesPerson.Range("A2").Value = "Chetan"
$salesPerson.Range("B2").Value = 0.1
$salesPerson.Range("B2").NumberFormat = "0.00%"
$salesPerson.Range("B2").HorizontalAlignment = -4108

$salesPerson.Range("A3").Value = "Rahul"
$salesPerson.Range("B3").Value = 0.1
$salesPerson.Range("B2").Copy()
$salesPerson.Range("B3").PasteSpecial(-4122)
$salesPerson.Range("B3").Value = 0.1
$excel.CutCopyMode = $false

[void]$salesPerson.Range("B5").Select()
$salesPerson.PageSetup.Orientation = 1

Write-Host "done"
